$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.790.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.078.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.01%  "

# Row 10
$ws.Range("E10").Value = "  +2.21%  "

# Row 11
$ws.Range("E11").Value = "  +2.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.370.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.083.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.673.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.65%  "

# Row 29
$ws.Range("E29").Value = "  +0.70%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "

# Row 32
$ws.Range("E32").Value = "  +3.63%  "

# Row 33
$ws.Range("E33").Value = "  +1.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "

# Row 35
$ws.Range("E35").Value = "  +0.65%  "

# Row 36
$ws.Range("E36").Value = "  +3.87%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.50%  "

# Row 41
$ws.Range("E41").Value = "  -0.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.56%  "

# Row 43
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.447.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "

# Row 45
$ws.Range("E45").Value = "  -0.50%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.03%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.92%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.267.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "

